$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.743.49"

$ws.Range("D3").Value = "3.609.43"
$ws.Range("E3").Value = "  +2.22%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "202.61"
$ws.Range("E5").Value = "  +3.37%  "

$ws.Range("D6").Value = "602.04"
$ws.Range("E6").Value = "  -0.59%  "

$ws.Range("E7").Value = "  +0.66%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("E9").Value = "  +6.40%  "

$ws.Range("D10").Value = "0.647"
$ws.Range("E10").Value = "  +0.11%  "

$ws.Range("D11").Value = "53.95"
$ws.Range("E11").Value = "  +0.91%  "

$ws.Range("E12").Value = "  +0.09%  "

$ws.Range("D13").Value = "9.62"
$ws.Range("E13").Value = "  +1.56%  "

$ws.Range("D14").Value = "4.178.39"
$ws.Range("E14").Value = "  +2.08%  "

$ws.Range("D15").Value = "677.15"
$ws.Range("E15").Value = "  +13.24%  "

$ws.Range("D16").Value = "70.795.77"
$ws.Range("E16").Value = "  +1.16%  "

$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").Value = "12.82"
$ws.Range("E17").Value = "  +0.79%  "

$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").Value = "19.18"
$ws.Range("E18").Value = "  +0.89%  "

$ws.Range("D19").Value = "3.610.69"
$ws.Range("E19").Value = "  +2.28%  "

$ws.Range("E20").Value = "  +0.47%  "

$ws.Range("E21").Value = "  +1.48%  "

$ws.Range("D22").Value = "18.74"
$ws.Range("E22").Value = "  +4.60%  "

$ws.Range("D23").Value = "109.99"
$ws.Range("E23").Value = "  +6.60%  "

$ws.Range("D24").Value = "5.36"
$ws.Range("E24").Value = "  +3.76%  "

$ws.Range("E25").Value = "  +0.30%  "

$ws.Range("E26").Value = "  -0.64%  "

$ws.Range("E27").Value = "  -1.84%  "

$ws.Range("E28").Value = "  -0.80%  "

$ws.Range("D29").Value = "10.15"
$ws.Range("E29").Value = "  +6.40%  "

$ws.Range("D30").Value = "34.39"
$ws.Range("E30").Value = "  +3.27%  "

$ws.Range("D31").Value = "4.47"
$ws.Range("E31").Value = "  +5.56%  "

$ws.Range("D32").Value = "7.18"
$ws.Range("E32").Value = "  +1.30%  "

$ws.Range("D33").Value = "12.28"
$ws.Range("E33").Value = "  -0.45%  "

$ws.Range("D35").Value = "63.59"
$ws.Range("E35").Value = "  +0.18%  "

$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D36").Value = "0.0₃0856"
$ws.Range("E36").Value = "  +5.44%  "

$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "3.894.13"
$ws.Range("E37").Value = "  +3.14%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  -0.03%  "

$ws.Range("D39").Value = "513.49"
$ws.Range("E39").Value = "  +1.05%  "

$ws.Range("D40").Value = "3.02"
$ws.Range("E40").Value = "  -4.37%  "

$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "3.61"
$ws.Range("E41").Value = "  +1.14%  "

$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").Value = "36.94"
$ws.Range("E42").Value = "  +1.31%  "

$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.140"
$ws.Range("E43").Value = "  +5.02%  "

$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").Value = "0.384"
$ws.Range("E44").Value = "  -1.76%  "

$ws.Range("E45").Value = "  +4.12%  "

$ws.Range("D46").Value = "3.06"
$ws.Range("E46").Value = "  +8.56%  "

$ws.Range("D47").Value = "3.41"
$ws.Range("E47").Value = "  +5.16%  "

$ws.Range("E48").Value = "  +1.65%  "

$ws.Range("D49").Value = "8.64"
$ws.Range("E49").Value = "  +1.94%  "

$ws.Range("E50").Value = "  -0.20%  "

$ws.Range("D51").Value = "1.81"
$ws.Range("E51").Value = "  +20.24%  "
